$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the obsolete "PIOMPRIDE 4/30 MG 30 TAB." row entirely (row 64).
# Deleting the row shifts every following row up by one and drops the
# now-unused shared-string entry automatically.
$ws.Rows.Item(64).Delete()

# Column A ("م") is a plain typed-in sequence number, not a formula, so
# deleting the row doesn't renumber it automatically. Re-sequence the
# rows that shifted up (old rows 65..93, now sitting at 64..92) so the
# numbering stays consecutive (...,60,61,62,... instead of skipping 61).
for ($r = 64; $r -le 92; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 3
}

# The running-balance total (previously row 94, now row 93 after the
# shift) was updated by hand to reflect the removed row's balance
# (5958.17 - 19.5 = 5938.67).
$ws.Cells.Item(93, 11).Value = 5938.67

# The footer row (previously row 95, now row 94) got a slightly taller
# row height in the new version of the file.
$ws.Rows.Item(94).RowHeight = 17.25
